# Updated cryptos list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.493.40"
$ws.Range("E2").Value = "  +0.34%  "

$ws.Range("D3").Value = "2.104.96"
$ws.Range("E3").Value = "  +4.44%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.17"
$ws.Range("E5").Value = "  +1.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5270"
$ws.Range("E7").Value = "  +2.70%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4393"
$ws.Range("E8").Value = "  +3.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08869"
$ws.Range("E9").Value = "  +1.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.64"
$ws.Range("E10").Value = "  +9.60%  "

$ws.Range("E11").Value = "  +2.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.69"
$ws.Range("E12").Value = "  +0.04%  "

$ws.Range("D13").Value = "2.110.49"
$ws.Range("E13").Value = "  +4.94%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.744"
$ws.Range("E14").Value = "  +2.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.771"
$ws.Range("E15").Value = "  +3.99%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.50"
$ws.Range("E16").Value = "  +2.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.004"
$ws.Range("E17").Value = "  +0.21%  "

$ws.Range("E18").Value = "  +1.45%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06643"
$ws.Range("E19").Value = "  +1.77%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.04"
$ws.Range("E20").Value = "  +0.66%  "

$ws.Range("E21").Value = "  +0.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.304"
$ws.Range("E22").Value = "  +1.37%  "

$ws.Range("D23").Value = "30.552.85"
$ws.Range("E23").Value = "  +0.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.30"
$ws.Range("E24").Value = "  +3.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.348"
$ws.Range("E25").Value = "  +3.39%  "

$ws.Range("D26").Value = "2.348.43"
$ws.Range("E26").Value = "  +4.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.46"
$ws.Range("E27").Value = "  -0.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.604"
$ws.Range("E28").Value = "  +6.88%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.66"
$ws.Range("E29").Value = "  -0.56%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.83"
$ws.Range("E30").Value = "  +1.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.212"
$ws.Range("E31").Value = "  +5.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1073"
$ws.Range("E32").Value = "  +1.86%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.683"
$ws.Range("E33").Value = "  +22.61%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.242"
$ws.Range("E34").Value = "  +1.98%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.930"
$ws.Range("E35").Value = "  +2.63%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.23"
$ws.Range("E36").Value = "  +11.67%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02586"
$ws.Range("E37").Value = "  +2.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.501"
$ws.Range("E38").Value = "  +0.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06715"
$ws.Range("E39").Value = "  +0.66%  "

$ws.Range("E40").Value = "  +2.42%  "

$ws.Range("E41").Value = "  +3.03%  "

$ws.Range("E42").Value = "  +2.74%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.268"
$ws.Range("E43").Value = "  +2.55%  "

$ws.Range("E44").Value = "  +0.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.04"
$ws.Range("E45").Value = "  +2.23%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6400"
$ws.Range("E46").Value = "  +3.47%  "

$ws.Range("E47").Value = "  +0.76%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.632"
$ws.Range("E48").Value = "  -1.04%  "

$ws.Range("E49").Value = "  -0.75%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.215"
$ws.Range("E50").Value = "  +10.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "82.34"
$ws.Range("E51").Value = "  +1.39%  "
